# Inventory Requisition to Post Payment API
# Add Project / Division columns to the ObjectName sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ObjectName")

$ws.Range("F1").Value = "Project"
$ws.Range("G1").Value = "Division"
$ws.Range("F2").Value = "100 Home Project"
$ws.Range("G2").Value = "Colorado (100)"

$ws.Columns.Item(7).ColumnWidth = 13.166666666666666

$ws.Range("F1:G1048576").Select()
